# Fix for drawing bitmap issue: add the missing time-log entry for row 85
# (a Coding session on 2014-10-14, 09:09 - 09:31) that had been left blank,
# and advance the active-cell selection to C86 - matching the recorded
# workbook state after this row was filled in.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A85").Value = 41926                  # Date       2014-10-14
$ws.Range("B85").Value = 0.38125000000000003    # Start Time 09:09 AM
$ws.Range("C85").Value = 0.39652777777777781    # Stop Time  09:31 AM
$ws.Range("D85").Value = 0                      # Interruption (mins)
$ws.Range("F85").Value = "Coding"                # Activity

# E85 already carries the shared formula (=...(C85-B85)*24-D85/60...) from
# the template row, so it recalculates to the correct Delta automatically.

$ws.Range("C86").Select()
